$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$rows = @{
    2  = @{D=44776; J=80;  K=20000; L=20000; M=20000; P=1333}
    3  = @{D=44769; J=50;  K=20000; L=20000; M=20000; P=1333}
    4  = @{D=44819; J=100; K=20000; L=20000; M=20000; P=1333}
    5  = @{D=44756; J=80;  K=20000; L=20000; M=20000; P=1333}
    6  = @{D=44757; J=30;  K=20000; L=20000; M=20000; P=1333}
    7  = @{D=44838; J=10;  K=20000; L=20000; M=20000; P=1333}
    8  = @{D=44508; J=40;  K=10000; L=10000; M=10000; P=667}
    9  = @{D=44845; J=20;  K=16000; L=16000; M=16000; P=1067}
    10 = @{D=44830; J=25;  K=12000; L=12000; M=12000; P=800}
    11 = @{D=44837; J=80;  K=16000; L=16000; M=16000; P=1067}
    12 = @{D=44771; J=40;  K=20000; L=20000; M=20000; P=1333}
    13 = @{D=44767; J=50;  K=20000; L=20000; M=20000; P=1333}
    14 = @{D=44839; J=80;  K=16000; L=16000; M=16000; P=1067}
    15 = @{D=44825; J=30;  K=20000; L=20000; M=20000; P=1333}
    16 = @{D=44755; J=50;  K=20000; L=20000; M=20000; P=1333}
    17 = @{D=44811; J=30;  K=20000; L=20000; M=20000; P=1333}
    18 = @{D=44826; J=50;  K=20000; L=20000; M=20000; P=1333}
    19 = @{D=44518; J=50;  K=10000; L=10000; M=10000; P=667}
    20 = @{D=44749; J=50;  K=20000; L=20000; M=20000; P=1333}
    21 = @{D=44525; J=40;  K=8000;  L=8000;  M=8000;  P=533}
    22 = @{D=44827; J=20;  K=20000; L=20000; M=20000; P=1333}
    23 = @{D=44841; J=20;  K=16000; L=16000; M=16000; P=1067}
    24 = @{D=44813; J=20;  K=20000; L=20000; M=20000; P=1333}
    25 = @{D=44824; J=20;  K=20000; L=20000; M=20000; P=1333}
    26 = @{D=44812; J=80;  K=20000; L=20000; M=20000; P=1333}
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
